# "dev of alliance region"
#
# Reworks the "normal" sheet's unit-type taxonomy from Chinese free-text
# descriptions to a small set of shared English category codes, and renames
# the archer_1/2/3 units to ranger_1/2/3 so "archer" is freed up to become
# the generic category name. Also switches the active sheet/selection back
# to "normal" (it was "special" before).

$wb = $excel.ActiveWorkbook
$wsNormal = $wb.Worksheets.Item("normal")

# --- Header row: rename STR_description -> STR_category, the INT_atk* set
#     -> INT_infantry/INT_archer/INT_cavalry/INT_siege/INT_wall, and
#     INT_upkeep -> INT_consume ---
$wsNormal.Range("B1").Value = "STR_category"
$wsNormal.Range("E1").Value = "INT_infantry"
$wsNormal.Range("F1").Value = "INT_archer"
$wsNormal.Range("G1").Value = "INT_cavalry"
$wsNormal.Range("H1").Value = "INT_siege"
$wsNormal.Range("I1").Value = "INT_wall"
$wsNormal.Range("M1").Value = "INT_consume"

# --- Unit name renames: archer_1/2/3 -> ranger_1/2/3 (rows 8-10) ---
$wsNormal.Range("A8").Value = "ranger_1"
$wsNormal.Range("A9").Value = "ranger_2"
$wsNormal.Range("A10").Value = "ranger_3"

# --- Column B (category) rewrite for every data row, by unit group ---
$wsNormal.Range("B2:B7").Value = "infantry"    # swordsman_*, sentinel_*
$wsNormal.Range("B8:B13").Value = "archer"     # ranger_* (ex archer_*), crossbowman_*
$wsNormal.Range("B14:B19").Value = "cavalry"   # lancer_*, horseArcher_*
$wsNormal.Range("B20:B25").Value = "siege"     # catapult_*, ballista_*

# --- Switch the active sheet/selection to "normal" (M2) ---
$wsNormal.Activate()
$wsNormal.Range("M2").Select() | Out-Null
